# The target adds row 1 / cell A1 on the sheet with the (string) value "3".
# Excel auto-coerces a numeric-looking literal like "3" into a Number when
# the cell is in the default "General" format, so the cell has to be primed
# as Text first to keep it stored as a string - exactly what typing a
# text-formatted "3" into a cell in the Excel UI would do.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("A1")
$cell.NumberFormat = "@"
$cell.Value = "3"
# Drop back to the workbook's default style now that the text value is
# committed, so we don't leave a stray text-format style behind on A1.
$cell.Style = "Normal"
